# Applies the "Diabolos Profits" data refresh committed by the scheduled
# runner. For each affected leve row (identified by sheet + row number),
# the currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ
# derived columns (H, I, J, K, L, M, N) are updated to their newly
# recalculated values.

$wb = $excel.ActiveWorkbook

function Set-Cells {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values
    )

    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# --- ALC -------------------------------------------------------------
Set-Cells "ALC" 69 @{
    H = 5037
    I = 2500
    J = 7574
    K = 7500
    L = 22722
    M = -6626
    N = -24470
}

Set-Cells "ALC" 72 @{
    H = 5037
    I = 2500
    J = 7574
    K = 22500
    L = 68166
    M = -18132
    N = -76902
}

Set-Cells "ALC" 129 @{
    H = 1294.7
    I = 690.25
    K = 2070.75
    M = 2929.25
}

Set-Cells "ALC" 131 @{
    H = 12921.357
    I = 3575
    K = 10725
    M = -5685
}

Set-Cells "ALC" 138 @{
    H = 1979.3829
    I = 1204.6
    K = 3613.8
    M = 1526.2
}

# --- ARM -------------------------------------------------------------
Set-Cells "ARM" 61 @{
    H = 15876483
    I = 22225376
    K = 22225376
    M = -22225164
}

Set-Cells "ARM" 74 @{
    H = 2575.8572
    I = 2450.0557
    K = 2450.0557
    M = -1576.0557
}

Set-Cells "ARM" 77 @{
    H = 2575.8572
    I = 2450.0557
    K = 12250.2785
    M = -7882.2785
}

Set-Cells "ARM" 136 @{
    H = 15876483
    I = 22225376
    K = 66676128
    M = -66673578
}

Set-Cells "ARM" 138 @{
    H = 69000
    J = 69000
    L = 69000
    N = -79280
}

# --- BSM -------------------------------------------------------------
Set-Cells "BSM" 20 @{
    H = 12055.179
    I = 14342
    J = 6338.125
    K = 14342
    L = 6338.125
    M = -14095
    N = -6832.125
}

Set-Cells "BSM" 37 @{
    H = 469.66666
    I = 465.875
    J = 500
    K = 465.875
    L = 500
    M = -328.875
    N = -774
}

Set-Cells "BSM" 134 @{
    H = 1288.5
    I = 1199.6875
    K = 3599.0625
    M = -1064.0625
}

# --- CRP -------------------------------------------------------------
Set-Cells "CRP" 31 @{
    H = 2664.141
    I = 1558.7428
    K = 1558.7428
    M = -1263.7428
}

Set-Cells "CRP" 34 @{
    H = 2664.141
    I = 1558.7428
    K = 1558.7428
    M = -1356.7428
}

Set-Cells "CRP" 86 @{
    H = 5959.3125
    I = 6033.25
    J = 5737.5
    K = 6033.25
    L = 5737.5
    M = -4910.25
    N = -7983.5
}

Set-Cells "CRP" 89 @{
    H = 5959.3125
    I = 6033.25
    J = 5737.5
    K = 30166.25
    L = 28687.5
    M = -24550.25
    N = -39919.5
}

Set-Cells "CRP" 105 @{
    H = 2473.45
    I = 1375.8
    J = 3571.1
    K = 1375.8
    L = 3571.1
    M = 371.2
    N = -7065.1
}

# --- CUL -------------------------------------------------------------
Set-Cells "CUL" 129 @{
    H = 2745.4
    I = 1741.5
    J = 3414.6667
    K = 5224.5
    L = 10244.0001
    M = -224.5
    N = -20244.0001
}

# --- GSM -------------------------------------------------------------
Set-Cells "GSM" 132 @{
    H = 3058.6216
    I = 2521.6775
    K = 7565.0325
    M = -5035.0325
}

Set-Cells "GSM" 136 @{
    H = 39825.75
    J = 39825.75
    L = 119477.25
    N = -124577.25
}

Set-Cells "GSM" 141 @{
    H = 82233
    I = 59700
    J = 93499.5
    K = 59700
    L = 93499.5
    M = -54520
    N = -103859.5
}

# --- LTW -------------------------------------------------------------
Set-Cells "LTW" 19 @{
    H = 2308.2
    I = 1000
    J = 2635.25
    K = 1000
    L = 2635.25
    M = -830
    N = -2975.25
}

Set-Cells "LTW" 22 @{
    H = 686.6
    J = 686.6
    L = 686.6
    N = -1276.6
}

Set-Cells "LTW" 27 @{
    H = 686.6
    J = 686.6
    L = 686.6
    N = -900.6
}

Set-Cells "LTW" 46 @{
    H = 1528.9166
    I = 909.5
    K = 909.5
    M = -721.5
}

Set-Cells "LTW" 100 @{
    H = 3049.25
    I = 2845.3635
    J = 3497.8
    K = 2845.3635
    L = 3497.8
    M = -2304.3635
    N = -4579.8
}

Set-Cells "LTW" 122 @{
    H = 3706.7778
    I = 2893.5
    J = 5333.3335
    K = 8680.5
    L = 16000.0005
    M = -6230.5
    N = -20900.0005
}

# --- WVR -------------------------------------------------------------
Set-Cells "WVR" 96 @{
    H = 1436.875
    I = 1542.1428
    K = 1542.1428
    M = -169.1428000000001
}

Set-Cells "WVR" 122 @{
    H = 3499
    I = 3499
    K = 10497
    M = -8047
}

Set-Cells "WVR" 126 @{
    H = 2291.3333
    I = 2049.6
    K = 6148.799999999999
    M = -3678.799999999999
}
